# This script reproduces the "Updated cryptos list ... with GitHub Actions"
# commit: it refreshes the Price (column D) and Volume(1h) (column E) figures
# for every coin row (2-51) on the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into column D (Price). Values are kept as plain text
# exactly like the existing cells (dotted thousands separators, fixed decimal
# places, etc.).
$priceUpdates = @{
    2 = "27.865.91"
    3 = "1.767.22"
    4 = "1.002"
    5 = "327.28"
    7 = "0.4473"
    8 = "0.3550"
    9 = "0.07439"
    10 = "42.04"
    12 = "1.002"
    13 = "20.83"
    14 = "6.031"
    15 = "7.200"
    16 = "1.772.52"
    17 = "92.90"
    19 = "0.06428"
    21 = "17.16"
    22 = "5.797"
    23 = "27.890.46"
    25 = "2.122"
    26 = "163.02"
    27 = "20.21"
    28 = "1.974.75"
    29 = "2.170"
    30 = "125.21"
    31 = "1.102"
    32 = "0.09151"
    33 = "5.576"
    34 = "3.636"
    35 = "11.84"
    36 = "0.02292"
    37 = "0.06100"
    38 = "0.2093"
    39 = "0.6320"
    40 = "4.968"
    41 = "1.182"
    42 = "1.391"
    43 = "7.943"
    44 = "13.23"
    45 = "3.734"
    46 = "0.5875"
    47 = "122.25"
    48 = "1.952"
    49 = "0.06904"
    50 = "1.138"
    51 = "72.91"
}

# New values to write into column E (Volume(1h)), including the original
# double-space padding around the percentage text.
$volumeUpdates = @{
    2 = "  +0.51%  "
    3 = "  +0.51%  "
    6 = "  +0.08%  "
    7 = "  -3.19%  "
    8 = "  -1.61%  "
    9 = "  -1.14%  "
    11 = "  -0.53%  "
    12 = "  +0.17%  "
    13 = "  -0.11%  "
    14 = "  -0.05%  "
    15 = "  +1.04%  "
    16 = "  +0.88%  "
    17 = "  +0.42%  "
    18 = "  -0.95%  "
    19 = "  +0.28%  "
    20 = "  +0.09%  "
    21 = "  +1.94%  "
    22 = "  -0.57%  "
    23 = "  +0.44%  "
    24 = "  +0.07%  "
    25 = "  +0.75%  "
    26 = "  -1.00%  "
    27 = "  -1.12%  "
    28 = "  +1.17%  "
    29 = "  +3.51%  "
    30 = "  -1.33%  "
    31 = "  +2.92%  "
    32 = "  -0.82%  "
    33 = "  +0.07%  "
    34 = "  -0.80%  "
    35 = "  -1.03%  "
    36 = "  -0.58%  "
    37 = "  +0.85%  "
    38 = "  -0.63%  "
    39 = "  -1.00%  "
    40 = "  -0.34%  "
    41 = "  -1.76%  "
    42 = "  +0.94%  "
    43 = "  +1.11%  "
    44 = "  -0.07%  "
    45 = "  +0.53%  "
    46 = "  -0.93%  "
    47 = "  -1.06%  "
    48 = "  -0.43%  "
    49 = "  +0.45%  "
    50 = "  -1.56%  "
    51 = "  +0.39%  "
}

# Rows whose new Price text would otherwise be auto-recognized by Excel as a
# number (e.g. "1.002", "0.3550"), which would silently convert the cell to a
# numeric value and drop significant trailing zeros / exact formatting. Force
# those specific cells to Text format before assigning the value, then restore
# the default "Normal" cell style so no stray formatting is left behind.
$forceTextRows = @(4, 5, 7, 8, 9, 10, 12, 13, 14, 15, 17, 19, 21, 22, 25, 26, 27, 29, 30, 31, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 45, 46, 47, 48, 49, 50, 51)

foreach ($row in $forceTextRows) {
    $ws.Range("D$row").NumberFormat = "@"
}

foreach ($row in $priceUpdates.Keys) {
    $ws.Range("D$row").Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Range("E$row").Value = $volumeUpdates[$row]
}

foreach ($row in $forceTextRows) {
    $ws.Range("D$row").Style = "Normal"
}
